$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row for the "Unknown" payment type with a numeric placeholder (0)
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "Unknown"

# Resize the table to include the newly added row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B8"))

# Update selection to match the recorded cursor position after the edit
$ws.Range("B11").Select()
